## ------------------------------------------------------------------------
## Applies the "last update (safe backup)" edit:
##  1. Before overwriting the "Hours" sheet with new data, copy its current
##     (old) contents into a brand-new "Sheet1" tab as a backup, together
##     with a small summary block in columns F:H and a "real thing!" title.
##  2. Overwrite "Hours" (A4:C10) with the new period/weight/season figures,
##     using the same "+24 hour" running formula pattern as before, and
##     drop the now-unused 11th row.
##  3. Leave the selections / scroll position on "Hours" and "Blad3" the way
##     the author left them when they saved.
## ------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$hours = $wb.Worksheets.Item("Hours")
$blad3 = $wb.Worksheets.Item("Blad3")

## --- 1. snapshot of the current "Hours" values, taken before they get
##        overwritten below (kept as literals so the backup sheet's numbers
##        round-trip byte-exactly instead of drifting through COM marshalling)
$oldA = 672,  3528,  4032,  5376,  5880,  6888,  7536,  8232
$oldB = 3.0659999999999998, 8.375, 4.742, 5.2770000000000001, 6.98, 5.8559999999999999, 9.6780000000000008, 8.1690000000000005
$oldC = 1, 2, 2, 3, 3, 4, 4, 1

## --- 2. build the backup sheet ("Sheet1") from that snapshot -------------
$lastIndex = $wb.Worksheets.Count
$afterSheet = $wb.Worksheets.Item($lastIndex)
$backup = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$backup.Name = "Sheet1"

$backup.Range("A1").Value = "real thing!"

$backup.Range("A2").Value = "First hour"
$backup.Range("A2").Font.Bold = $true
$backup.Range("B2").Value = "Weight"
$backup.Range("B2").Font.Bold = $true

for ($i = 0; $i -lt 8; $i++) {
    $row = 3 + $i
    $backup.Cells.Item($row, 1).Value = $oldA[$i]
    $backup.Cells.Item($row, 2).Value = $oldB[$i]
    $backup.Cells.Item($row, 3).Value = $oldC[$i]
}

## small summary block tucked in columns F:H of the backup sheet
$backup.Range("F3").Value = 3528
$backup.Range("G3").Value = 13.035714285714301
$backup.Range("H3").Value = 2

$backup.Range("F4").Value = 5880
$backup.Range("G4").Value = 13.035714285714286
$backup.Range("H4").Value = 3

$backup.Range("F5").Value = 7536
$backup.Range("G5").Value = 13.035714285714286
$backup.Range("H5").Value = 4

$backup.Range("F6").Value = 8232
$backup.Range("G6").Value = 13.035714285714286
$backup.Range("H6").Value = 1

$backup.Range("A8").Select() | Out-Null

## --- 3. overwrite "Hours" with the new figures ----------------------------
$hours.Range("A4").Value = 7536
$hours.Range("B4").Value = 52.142857142857146
$hours.Range("C4").Value = 4

$hours.Range("A5").Formula = "=A4+24"
$hours.Range("B5").Value = 52.142857142857146
$hours.Range("C5").Value = 4

$hours.Range("A6:A10").Formula = "=A5+24"
$hours.Range("B6").Value = 52.142857142857146
$hours.Range("C6").Value = 4
$hours.Range("B7").Value = 52.142857142857146
$hours.Range("C7").Value = 4
$hours.Range("B8").Value = 52.142857142857146
$hours.Range("C8").Value = 4
$hours.Range("B9").Value = 52.142857142857146
$hours.Range("C9").Value = 4
$hours.Range("B10").Value = 52.142857142857146
$hours.Range("C10").Value = 4

## row 11 no longer exists
$hours.Range("A11:C11").ClearContents() | Out-Null

## --- 4. restore the sheet selections -------------------------------------
$blad3.Range("R4").Select() | Out-Null
$hours.Range("F9").Select() | Out-Null
